# Updates cryptocurrency price/volume figures in the "cryptos" worksheet.
# Values that are plain numeric strings (e.g. "596.18") are written with a
# leading quote so Excel keeps them as literal text (quote-prefixed), matching
# the source data which stores prices/volumes as text, not numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{cell="D2"; value='68.009.49'; forceText=$false},
    @{cell="E2"; value='  +1.65%  '; forceText=$false},
    @{cell="D3"; value='2.623.66'; forceText=$false},
    @{cell="E3"; value='  +1.32%  '; forceText=$false},
    @{cell="E4"; value='  +0.13%  '; forceText=$false},
    @{cell="D5"; value='596.18'; forceText=$true},
    @{cell="E5"; value='  +0.70%  '; forceText=$false},
    @{cell="D6"; value='155.42'; forceText=$true},
    @{cell="E6"; value='  +0.08%  '; forceText=$false},
    @{cell="E7"; value='  +0.06%  '; forceText=$false},
    @{cell="E8"; value='  +1.11%  '; forceText=$false},
    @{cell="D9"; value='2.622.04'; forceText=$false},
    @{cell="E9"; value='  +1.20%  '; forceText=$false},
    @{cell="E10"; value='  +10.68%  '; forceText=$false},
    @{cell="E11"; value='  +1.02%  '; forceText=$false},
    @{cell="E12"; value='  +0.13%  '; forceText=$false},
    @{cell="D13"; value='0.356'; forceText=$true},
    @{cell="E13"; value='  -1.46%  '; forceText=$false},
    @{cell="D14"; value='27.77'; forceText=$true},
    @{cell="E14"; value='  -2.21%  '; forceText=$false},
    @{cell="E15"; value='  +3.35%  '; forceText=$false},
    @{cell="D16"; value='3.103.95'; forceText=$false},
    @{cell="E16"; value='  +1.37%  '; forceText=$false},
    @{cell="D17"; value='67.714.55'; forceText=$false},
    @{cell="E17"; value='  +1.56%  '; forceText=$false},
    @{cell="D18"; value='2.623.98'; forceText=$false},
    @{cell="E18"; value='  +1.54%  '; forceText=$false},
    @{cell="D19"; value='369.95'; forceText=$true},
    @{cell="E19"; value='  +3.58%  '; forceText=$false},
    @{cell="D20"; value='11.19'; forceText=$true},
    @{cell="E20"; value='  -0.97%  '; forceText=$false},
    @{cell="D21"; value='7.66'; forceText=$true},
    @{cell="E21"; value='  -2.09%  '; forceText=$false},
    @{cell="E22"; value='  -0.70%  '; forceText=$false},
    @{cell="D23"; value='2.04'; forceText=$true},
    @{cell="E23"; value='  -1.42%  '; forceText=$false},
    @{cell="E24"; value='  -0.03%  '; forceText=$false},
    @{cell="D25"; value='9.91'; forceText=$true},
    @{cell="E25"; value='  -6.75%  '; forceText=$false},
    @{cell="D26"; value='67.42'; forceText=$true},
    @{cell="E26"; value='  +0.47%  '; forceText=$false},
    @{cell="E27"; value='  +0.65%  '; forceText=$false},
    @{cell="D29"; value='576.38'; forceText=$true},
    @{cell="E29"; value='  -6.28%  '; forceText=$false},
    @{cell="E30"; value='  +0.16%  '; forceText=$false},
    @{cell="D31"; value='1.44'; forceText=$true},
    @{cell="E31"; value='  -2.28%  '; forceText=$false},
    @{cell="E32"; value='  -1.03%  '; forceText=$false},
    @{cell="E33"; value='  +0.33%  '; forceText=$false},
    @{cell="E34"; value='  -1.57%  '; forceText=$false},
    @{cell="E35"; value='  +0.02%  '; forceText=$false},
    @{cell="E36"; value='  -3.44%  '; forceText=$false},
    @{cell="D37"; value='4.95'; forceText=$true},
    @{cell="E37"; value='  -2.51%  '; forceText=$false},
    @{cell="D38"; value='158.85'; forceText=$true},
    @{cell="E38"; value='  +2.88%  '; forceText=$false},
    @{cell="D39"; value='19.39'; forceText=$true},
    @{cell="E39"; value='  +0.44%  '; forceText=$false},
    @{cell="E40"; value='  -0.37%  '; forceText=$false},
    @{cell="D41"; value='5.35'; forceText=$true},
    @{cell="E41"; value='  -3.24%  '; forceText=$false},
    @{cell="E42"; value='  +1.78%  '; forceText=$false},
    @{cell="E43"; value='  -4.13%  '; forceText=$false},
    @{cell="D44"; value='41.22'; forceText=$true},
    @{cell="E44"; value='  -0.93%  '; forceText=$false},
    @{cell="D45"; value='0.999'; forceText=$true},
    @{cell="E45"; value='  +0.00%  '; forceText=$false},
    @{cell="D46"; value='16.43'; forceText=$true},
    @{cell="E46"; value='  -0.08%  '; forceText=$false},
    @{cell="D47"; value='156.19'; forceText=$true},
    @{cell="E47"; value='  +0.11%  '; forceText=$false},
    @{cell="E48"; value='  -6.84%  '; forceText=$false},
    @{cell="D49"; value='3.76'; forceText=$true},
    @{cell="E49"; value='  -0.54%  '; forceText=$false},
    @{cell="E50"; value='  +2.56%  '; forceText=$false},
    @{cell="D51"; value='0.0543'; forceText=$true},
    @{cell="E51"; value='  -4.16%  '; forceText=$false}
)

foreach ($u in $updates) {
    if ($u.forceText) {
        # Leading apostrophe = Excel's quote-prefix, forces text storage
        $ws.Range($u.cell).Value = "'" + $u.value
    } else {
        $ws.Range($u.cell).Value = $u.value
    }
}

Write-Host "Updated $($updates.Count) cells on $($ws.Name)"
